$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are numeric-looking strings (e.g. "1.00",
# "0.0781") that must stay stored as text, matching the workbooks
# existing text-based Price column. Force a Text number format on just
# those specific cells before assigning, so Excel does not coerce them
# into numbers (which would drop formatting like trailing zeros).
$textCells = @("D4", "D5", "D6", "D8", "D10", "D11", "D12", "D16", "D22", "D23", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D38", "D39", "D40", "D44", "D46", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "46.755.98"
$ws.Range("E2").Value = "  +4.34%  "
$ws.Range("D3").Value = "2.261.94"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "300.69"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "100.26"
$ws.Range("E6").Value = "  +6.61%  "
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "35.54"
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "2.607.00"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "2.265.80"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "13.57"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "46.760.83"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").Value = "65.10"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "248.90"
$ws.Range("E23").Value = "  +4.70%  "
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "42.53"
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "19.82"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "2.77"
$ws.Range("E31").Value = "  +8.56%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "145.11"
$ws.Range("E32").Value = "  -4.66%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.40"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "3.18"
$ws.Range("E35").Value = "  +8.31%  "
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("D38").Value = "16.09"
$ws.Range("E38").Value = "  +18.76%  "
$ws.Range("D39").Value = "1.70"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("E41").Value = "  -3.84%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "1.96"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "1.790.87"
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("D46").Value = "90.94"
$ws.Range("E46").Value = "  +19.64%  "
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "71.45"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.485.84"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "93.52"
$ws.Range("E51").Value = "  -1.99%  "
